# Update scripts with new TPM values (Gnai2-Ccr5, YoungD0, lrc2p)
#
# The recalculated NATMI run dropped the "MuSCs" sending-cluster rows
# (old rows 8-10) entirely and refreshed the expression / specificity
# numbers for the remaining Sending/Target cluster combinations
# (ECs x {ECs,FAPs}, FAPs x {ECs,FAPs}, MuSCs x {ECs,FAPs}).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 8, 9 and 10 (MuSCs sending-cluster block) first so the
# remaining writes land on a sheet that is already the right shape.
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

function Set-RowValues {
    param($RowNum, $Values)
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $ws.Range($cols[$i] + $RowNum).Value = $Values[$i]
    }
}

# Row 2: ECs -> ECs
Set-RowValues 2 @(
    "ECs", "Gnai2", "Ccr5", "ECs", 3, 1,
    63.91118233333333, 191.733547,
    0.4067926910433548, 0.4067926910433549,
    1, 0.3333333333333333,
    0.001937666666666667, 0.005813,
    0.0230007399171451, 0.02300073991714511,
    0.1238385676345555, 1.114547108711,
    0.009356532886883765, 0.009356532886883769
)

# Row 3: ECs -> FAPs
Set-RowValues 3 @(
    "ECs", "Gnai2", "Ccr5", "FAPs", 3, 1,
    63.91118233333333, 191.733547,
    0.4067926910433548, 0.4067926910433549,
    1, 0.3333333333333333,
    0.082306, 0.246918,
    0.9769992600828549, 0.976999260082855,
    5.260273773127333, 47.342463958146,
    0.397436158156471, 0.3974361581564712
)

# Row 4: FAPs -> ECs
Set-RowValues 4 @(
    "FAPs", "Gnai2", "Ccr5", "ECs", 3, 1,
    57.4434, 172.3302,
    0.3656254573230189, 0.365625457323019,
    1, 0.3333333333333333,
    0.001937666666666667, 0.005813,
    0.0230007399171451, 0.02300073991714511,
    0.1113061614, 1.0017554526,
    0.008409656050973994, 0.008409656050973997
)

# Row 5: FAPs -> FAPs
Set-RowValues 5 @(
    "FAPs", "Gnai2", "Ccr5", "FAPs", 3, 1,
    57.4434, 172.3302,
    0.3656254573230189, 0.365625457323019,
    1, 0.3333333333333333,
    0.082306, 0.246918,
    0.9769992600828549, 0.976999260082855,
    4.7279364804, 42.5514283236,
    0.3572158012720449, 0.357215801272045
)

# Row 6: MuSCs -> ECs
Set-RowValues 6 @(
    "MuSCs", "Gnai2", "Ccr5", "ECs", 3, 1,
    35.755375, 107.266125,
    0.2275818516336261, 0.2275818516336262,
    1, 0.3333333333333333,
    0.001937666666666667, 0.005813,
    0.0230007399171451, 0.02300073991714511,
    0.06928199829166666, 0.623537984625,
    0.005234550979287339, 0.00523455097928734
)

# Row 7: MuSCs -> FAPs
Set-RowValues 7 @(
    "MuSCs", "Gnai2", "Ccr5", "FAPs", 3, 1,
    35.755375, 107.266125,
    0.2275818516336261, 0.2275818516336262,
    1, 0.3333333333333333,
    0.082306, 0.246918,
    0.9769992600828549, 0.976999260082855,
    2.94288189475, 26.48593705275,
    0.2223473006543388, 0.2223473006543388
)
